$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49 - new weekly task entry for WEEK 6 ("Write proposal")
$ws.Cells.Item(49, 1).Value = 39
$ws.Cells.Item(49, 2).Value = "W001"
$ws.Cells.Item(49, 3).Value = "Write document"
$ws.Cells.Item(49, 4).Value = "Write proposal"
$ws.Cells.Item(49, 5).Value = "Thanh"
$ws.Cells.Item(49, 6).Value = 3
$ws.Cells.Item(49, 7).Value = 3
$ws.Cells.Item(49, 8).Value = 1
$ws.Cells.Item(49, 9).Value = "Done"
$ws.Cells.Item(49, 10).Value = 42725
$ws.Cells.Item(49, 11).Value = 42725
$ws.Cells.Item(49, 12).Value = 42725
$ws.Cells.Item(49, 13).Value = 42725
$ws.Cells.Item(49, 14).Value = 42725
$ws.Cells.Item(49, 15).Value = "Yes"

# Update the view to match the weekly-update scroll position + current selection
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("A50").Select() | Out-Null
